$wb = $excel.ActiveWorkbook

# --- Sheet "cf": wind capacity factor column updates ---
$ws = $wb.Worksheets.Item("cf")
$ws.Range("B2").Value = 0.6
$ws.Range("B3").Value = 0.4
$ws.Range("B4").Value = 0.6
$ws.Range("B5").Value = 0.4
$ws.Range("B6").Value = 0.5

# --- Sheet "price": penalty + eBprice ---
$ws = $wb.Worksheets.Item("price")
$ws.Range("A2").Value = 5
$ws.Range("F2").Value = 0.004

# --- Sheet "electrolyzer": sInvcost, sOmc, sInv ---
$ws = $wb.Worksheets.Item("electrolyzer")
$ws.Range("E2").Value = 0.07000000000000001
$ws.Range("F2").Value = 0.0005
$ws.Range("I2").Value = 0.009065320247581964

# --- Sheet "hytank": update row 2, add Tank2 row 3 ---
$ws = $wb.Worksheets.Item("hytank")
$ws.Range("I2").Value = 0.05
$ws.Range("J2").Value = 0.0005
$ws.Range("K2").Value = 0.001
$ws.Range("N2").Value = 0.004012129359534564

$ws.Range("A3").Value = "Tank2"
$ws.Range("B3").Value = 500
$ws.Range("C3").Value = 100
$ws.Range("D3").Value = 100
$ws.Range("E3").Value = 30
$ws.Range("F3").Value = 20
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0.03
$ws.Range("J3").Value = 0.0005
$ws.Range("K3").Value = 0.001
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.08024258719069129
$ws.Range("N3").Value = 0.002407277615720738

# --- Sheet "fuel-cell": sInvcost, sOmc, sInv ---
$ws = $wb.Worksheets.Item("fuel-cell")
$ws.Range("E2").Value = 0.008
$ws.Range("F2").Value = 0.0005
$ws.Range("J2").Value = 0.0006419406975255303
